$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.09298733333333332
$ws.Range("H2").Value = 0.278962
$ws.Range("I2").Value = 0.009366289255637828
$ws.Range("J2").Value = 0.00936628925563783
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 8.509392159811776
$ws.Range("R2").Value = 76.58452943830599
$ws.Range("S2").Value = 0.00907173802297255
$ws.Range("T2").Value = 0.009071738022972552
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.09298733333333332
$ws.Range("H3").Value = 0.278962
$ws.Range("I3").Value = 0.009366289255637828
$ws.Range("J3").Value = 0.00936628925563783
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 0.012526509648
$ws.Range("R3").Value = 0.112738586832
$ws.Range("S3").Value = 0.0000133543279866194
$ws.Range("T3").Value = 0.0000133543279866194
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.09298733333333332
$ws.Range("H4").Value = 0.278962
$ws.Range("I4").Value = 0.009366289255637828
$ws.Range("J4").Value = 0.00936628925563783
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 0.2637658550077777
$ws.Range("R4").Value = 2.373892695069999
$ws.Range("S4").Value = 0.0002811969046786592
$ws.Range("T4").Value = 0.0002811969046786592
$ws.Range("I5").Value = 0.5157656456804108
$ws.Range("J5").Value = 0.5157656456804109
$ws.Range("M5").Value = 91.51130433333333
$ws.Range("N5").Value = 274.533913
$ws.Range("O5").Value = 0.9685519820468944
$ws.Range("P5").Value = 0.9685519820468945
$ws.Range("Q5").Value = 468.5796073414426
$ws.Range("R5").Value = 4217.216466072984
$ws.Range("S5").Value = 0.4995458383954582
$ws.Range("T5").Value = 0.4995458383954584
$ws.Range("I6").Value = 0.5157656456804108
$ws.Range("J6").Value = 0.5157656456804109
$ws.Range("O6").Value = 0.001425786415744213
$ws.Range("P6").Value = 0.001425786415744214
$ws.Range("S6").Value = 0.0007353716513186729
$ws.Range("T6").Value = 0.0007353716513186732
$ws.Range("I7").Value = 0.5157656456804108
$ws.Range("J7").Value = 0.5157656456804109
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.03002223153736139
$ws.Range("P7").Value = 0.03002223153736139
$ws.Range("S7").Value = 0.01548443563363399
$ws.Range("T7").Value = 0.015484435633634
$ws.Range("I8").Value = 0.4748680650639512
$ws.Range("J8").Value = 0.4748680650639513
$ws.Range("M8").Value = 91.51130433333333
$ws.Range("N8").Value = 274.533913
$ws.Range("O8").Value = 0.9685519820468944
$ws.Range("P8").Value = 0.9685519820468945
$ws.Range("Q8").Value = 431.4236384881967
$ws.Range("R8").Value = 3882.81274639377
$ws.Range("S8").Value = 0.4599344056284636
$ws.Range("T8").Value = 0.4599344056284637
$ws.Range("I9").Value = 0.4748680650639512
$ws.Range("J9").Value = 0.4748680650639513
$ws.Range("O9").Value = 0.001425786415744213
$ws.Range("P9").Value = 0.001425786415744214
$ws.Range("S9").Value = 0.000677060436438921
$ws.Range("T9").Value = 0.0006770604364389212
$ws.Range("I10").Value = 0.4748680650639512
$ws.Range("J10").Value = 0.4748680650639513
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509734999999999
$ws.Range("O10").Value = 0.03002223153736139
$ws.Range("P10").Value = 0.03002223153736139
$ws.Range("S10").Value = 0.01425659899904874
$ws.Range("T10").Value = 0.01425659899904874
